# Actualizacion automatica: inserta un nuevo cliente "AREVALO SAQUICELA LUIS
# MARCELO" (con valores en cero) antes de "CARRION ALVAREZ MARIO ANDRES" en
# ambas hojas, desplazando el resto de filas hacia abajo y actualizando el
# contador "de 8" -> "de 9" en la fila de totales de la primera hoja.

$wb = $excel.ActiveWorkbook

# --- Hoja 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Inserta una fila nueva en la posicion 4 (antes de CARRION ALVAREZ MARIO ANDRES),
# desplazando las filas existentes (4..10) hacia abajo (4..11 conservando valores).
$ws1.Rows.Item(4).EntireRow.Insert()

$ws1.Cells.Item(4, 1).Value = "VACA PANCHI CAROLINA"
$ws1.Cells.Item(4, 2).Value = "AREVALO SAQUICELA LUIS MARCELO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(4, $col).Value = 0
}

# La antigua fila de totales (fila 10) ahora es la fila 11; actualiza el
# texto "de 8" a "de 9" ya que ahora hay 9 clientes.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(11, $col)
    $cell.Value = ($cell.Value() -replace "de 8", "de 9")
}

# --- Hoja 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Misma insercion de fila para mantener ambas hojas sincronizadas.
$ws2.Rows.Item(4).EntireRow.Insert()

$ws2.Cells.Item(4, 1).Value = "VACA PANCHI CAROLINA"
$ws2.Cells.Item(4, 2).Value = "AREVALO SAQUICELA LUIS MARCELO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(4, $col).Value = 0
}
